$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.005.36"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "3.778.06"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "627.78"
$ws.Range("E5").Value = "  +4.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.29"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").Value = "3.778.98"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.458"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.79"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.29"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "4.411.40"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "3.774.38"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "69.052.71"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.62"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.21"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.53"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.707"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.99"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.01"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("E27").Value = "  +3.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "3.926.49"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.74"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("E35").Value = "  +18.24%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "3.727.66"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.93"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.31"
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.967"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "154.97"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.06"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  +3.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.37"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.35"
$ws.Range("E51").Value = "  -1.70%  "
